$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 105
$ws.Range("B105").Value = 6077497
$ws.Range("F105").Value = "Deportes Copiapo"
$ws.Range("G105").Value = "Nublense"
$ws.Range("H105").Value = 1
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = "D"
$ws.Range("K105").Value = 2.6
$ws.Range("L105").Value = 3.4
$ws.Range("M105").Value = 2.6
$ws.Range("N105").Value = 2.8
$ws.Range("O105").Value = 3.2
$ws.Range("P105").Value = 2.7
$ws.Range("T105").Value = 2.25
$ws.Range("U105").Value = 2
$ws.Range("V105").Value = 1.85
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = 2.2
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = 0
$ws.Range("AA105").Value = -0
$ws.Range("AB105").Value = -0.5
$ws.Range("AC105").Value = 0.425

# Row 106
$ws.Range("B106").Value = 6077763
$ws.Range("F106").Value = "Curico Unido"
$ws.Range("G106").Value = "Magallanes"
$ws.Range("H106").Value = 3
$ws.Range("I106").Value = 4
$ws.Range("J106").Value = "A"
$ws.Range("K106").Value = 2.15
$ws.Range("L106").Value = 3.5
$ws.Range("M106").Value = 3.2
$ws.Range("N106").Value = 2.625
$ws.Range("O106").Value = 3.5
$ws.Range("P106").Value = 2.6
$ws.Range("T106").Value = 2.75
$ws.Range("U106").Value = 1.975
$ws.Range("V106").Value = 1.875
$ws.Range("W106").Value = -1
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = 1.6
$ws.Range("Z106").Value = -1
$ws.Range("AA106").Value = 0.8999999999999999
$ws.Range("AB106").Value = 0.9750000000000001
$ws.Range("AC106").Value = -1

# Row 116
$ws.Range("B116").Value = 6078267
$ws.Range("F116").Value = "Huachipato"
$ws.Range("G116").Value = "Audax Italiano"
$ws.Range("H116").Value = 2
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = "H"
$ws.Range("K116").Value = 1.5
$ws.Range("L116").Value = 4.333
$ws.Range("M116").Value = 6
$ws.Range("N116").Value = 1.444
$ws.Range("O116").Value = 4.75
$ws.Range("P116").Value = 7
$ws.Range("Q116").Value = -1.25
$ws.Range("R116").Value = 2.025
$ws.Range("S116").Value = 1.825
$ws.Range("T116").Value = 2.75
$ws.Range("U116").Value = 1.8
$ws.Range("V116").Value = 2.05
$ws.Range("W116").Value = 0.444
$ws.Range("X116").Value = -1
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = 1.025
$ws.Range("AC116").Value = 1.05

# Row 117
$ws.Range("B117").Value = 6143704
$ws.Range("F117").Value = "Curico Unido"
$ws.Range("G117").Value = "Colo Colo"
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 1
$ws.Range("J117").Value = "A"
$ws.Range("K117").Value = 6.5
$ws.Range("L117").Value = 4.75
$ws.Range("M117").Value = 1.4
$ws.Range("N117").Value = 12
$ws.Range("O117").Value = 8.5
$ws.Range("P117").Value = 1.166
$ws.Range("Q117").Value = 2
$ws.Range("R117").Value = 2
$ws.Range("S117").Value = 1.8
$ws.Range("T117").Value = 3.25
$ws.Range("U117").Value = 1.875
$ws.Range("V117").Value = 1.925
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = -1
$ws.Range("Y117").Value = 0.1659999999999999
$ws.Range("Z117").Value = 1
$ws.Range("AC117").Value = 0.925

# Row 118
$ws.Range("B118").Value = 6078268
$ws.Range("F118").Value = "OHiggins"
$ws.Range("G118").Value = "Palestino"
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = "A"
$ws.Range("K118").Value = 3.1
$ws.Range("L118").Value = 3.3
$ws.Range("M118").Value = 2.3
$ws.Range("N118").Value = 2.9
$ws.Range("O118").Value = 3.4
$ws.Range("P118").Value = 2.375
$ws.Range("Q118").Value = 0.25
$ws.Range("U118").Value = 2
$ws.Range("V118").Value = 1.8
$ws.Range("Y118").Value = 1.375
$ws.Range("AB118").Value = -1
$ws.Range("AC118").Value = 0.8

# Row 119
$ws.Range("B119").Value = 6077768
$ws.Range("F119").Value = "Union La Calera"
$ws.Range("G119").Value = "Universidad Catolica"
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 3
$ws.Range("J119").Value = "A"
$ws.Range("K119").Value = 2.05
$ws.Range("L119").Value = 3.5
$ws.Range("M119").Value = 3.4
$ws.Range("N119").Value = 2.05
$ws.Range("O119").Value = 3.6
$ws.Range("P119").Value = 3.4
$ws.Range("Q119").Value = -0.25
$ws.Range("R119").Value = 1.8
$ws.Range("S119").Value = 2
$ws.Range("T119").Value = 2.75
$ws.Range("U119").Value = 1.975
$ws.Range("V119").Value = 1.825
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = 2.4
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = 1
$ws.Range("AB119").Value = 0.4875
$ws.Range("AC119").Value = -0.5

# Row 120
$ws.Range("B120").Value = 6077499
$ws.Range("F120").Value = "Deportes Copiapo"
$ws.Range("G120").Value = "Everton de Vina"
$ws.Range("H120").Value = 2
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = "H"
$ws.Range("K120").Value = 2.1
$ws.Range("L120").Value = 3.4
$ws.Range("M120").Value = 3.4
$ws.Range("P120").Value = 2.4
$ws.Range("Q120").Value = 0.25
$ws.Range("R120").Value = 1.775
$ws.Range("S120").Value = 2.1
$ws.Range("T120").Value = 2.75
$ws.Range("U120").Value = 1.85
$ws.Range("V120").Value = 2
$ws.Range("W120").Value = 1.9
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 0.7749999999999999
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = -1
$ws.Range("AC120").Value = 1

# Row 121
$ws.Range("B121").Value = 6078998
$ws.Range("F121").Value = "Magallanes"
$ws.Range("G121").Value = "Coquimbo Unido"
$ws.Range("H121").Value = 2
$ws.Range("I121").Value = 3
$ws.Range("J121").Value = "A"
$ws.Range("K121").Value = 1.909
$ws.Range("L121").Value = 3.6
$ws.Range("M121").Value = 3.8
$ws.Range("N121").Value = 2.15
$ws.Range("O121").Value = 3.75
$ws.Range("P121").Value = 3.1
$ws.Range("Q121").Value = -0.25
$ws.Range("R121").Value = 1.85
$ws.Range("S121").Value = 1.95
$ws.Range("T121").Value = 3
$ws.Range("U121").Value = 1.85
$ws.Range("V121").Value = 1.95
$ws.Range("W121").Value = -1
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = 2.1
$ws.Range("Z121").Value = -1
$ws.Range("AA121").Value = 0.95
$ws.Range("AB121").Value = 0.8500000000000001

# Row 122
$ws.Range("B122").Value = 6078269
$ws.Range("F122").Value = "Universidad de Chile"
$ws.Range("G122").Value = "Nublense"
$ws.Range("H122").Value = 3
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = "H"
$ws.Range("K122").Value = 1.85
$ws.Range("L122").Value = 3.4
$ws.Range("M122").Value = 4.333
$ws.Range("N122").Value = 1.8
$ws.Range("O122").Value = 3.6
$ws.Range("P122").Value = 4.5
$ws.Range("Q122").Value = -0.75
$ws.Range("R122").Value = 1.925
$ws.Range("S122").Value = 1.925
$ws.Range("T122").Value = 2.5
$ws.Range("U122").Value = 2.025
$ws.Range("V122").Value = 1.825
$ws.Range("W122").Value = 0.8
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = 0.925
$ws.Range("AA122").Value = -1
$ws.Range("AB122").Value = 1.025

# Row 136
$ws.Range("G136").Value = "Cobresal"
$ws.Range("H136").Value = 2
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = "H"
$ws.Range("K136").Value = 2.625
$ws.Range("N136").Value = 2.2
$ws.Range("O136").Value = 3.4
$ws.Range("P136").Value = 3.25
$ws.Range("Q136").Value = -0.25
$ws.Range("R136").Value = 1.9
$ws.Range("S136").Value = 1.9
$ws.Range("T136").Value = 2.5
$ws.Range("U136").Value = 1.9
$ws.Range("V136").Value = 1.9
$ws.Range("W136").Value = 1.2
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = 0.8999999999999999
$ws.Range("AA136").Value = -1
$ws.Range("AB136").Value = 0.8999999999999999
$ws.Range("AC136").Value = -1

# Row 137
$ws.Range("B137").Value = 7723528
$ws.Range("F137").Value = "Palestino"
$ws.Range("G137").Value = "Universidad Catolica"
$ws.Range("K137").Value = 1.95
$ws.Range("L137").Value = 3.6
$ws.Range("M137").Value = 3.4
$ws.Range("N137").Value = 2.2
$ws.Range("O137").Value = 3.6
$ws.Range("P137").Value = 3.2
$ws.Range("Q137").Value = -0.25
$ws.Range("R137").Value = 1.925
$ws.Range("S137").Value = 1.925
$ws.Range("T137").Value = 2.75
$ws.Range("U137").Value = 2.025
$ws.Range("V137").Value = 1.825

# Row 138
$ws.Range("B138").Value = 7723533
$ws.Range("F138").Value = "OHiggins"
$ws.Range("G138").Value = "Everton de Vina"
$ws.Range("K138").Value = 3
$ws.Range("L138").Value = 3.2
$ws.Range("M138").Value = 2.375
$ws.Range("N138").Value = 2.75
$ws.Range("O138").Value = 3.1
$ws.Range("P138").Value = 2.75
$ws.Range("Q138").Value = 0
$ws.Range("R138").Value = 1.925
$ws.Range("S138").Value = 1.925
$ws.Range("T138").Value = 2.25
$ws.Range("U138").Value = 2
$ws.Range("V138").Value = 1.85

# Row 139
$ws.Range("N139").Value = 1.95
$ws.Range("O139").Value = 3.6
$ws.Range("P139").Value = 3.8
$ws.Range("Q139").Value = -0.5
$ws.Range("R139").Value = 1.975
$ws.Range("S139").Value = 1.875

# Row 140
$ws.Range("N140").Value = 2.15
$ws.Range("O140").Value = 3.2
$ws.Range("P140").Value = 3.75
$ws.Range("Q140").Value = -0.25
$ws.Range("R140").Value = 1.8
$ws.Range("S140").Value = 2.05

# Row 141
$ws.Range("P141").Value = 5.25
$ws.Range("Q141").Value = -0.75
$ws.Range("R141").Value = 1.875
$ws.Range("S141").Value = 1.975

# Row 143
$ws.Range("N143").Value = 3.6
$ws.Range("O143").Value = 3.4
$ws.Range("P143").Value = 2.1
$ws.Range("Q143").Value = 0.5
$ws.Range("R143").Value = 1.8
$ws.Range("S143").Value = 2.05
